$d = $word.ActiveDocument

$hdr = $d.Sections(1).Headers(2)
$shp = $hdr.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

Write-Output "Before: [$($tr.Text)]"
$tr.Text = "TESTVALUE"
Write-Output "After: [$($tr.Text)]"
